# Timesheet update: fill in hours for "Development" and "Testing" rows
# (row 2 was "Analyse/onderzoek", becomes "Development"; row 3 was
# "Development", becomes "Testing"), touch up the last-day-of-month
# column formatting on every task row, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# --- Row 2: task changes to "Development", hours filled in -----------------
$ws.Range("A2").Value = "Development"
$ws.Range("D2").Value = 8
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 8
$ws.Range("G2").Value = 8
$ws.Range("L2").Value = 8
$ws.Range("M2").Value = 8
$ws.Range("N2").Value = 8
$ws.Range("R2").Value = 8
$ws.Range("S2").Value = 8
$ws.Range("T2").Value = 4

# --- Row 3: task changes to "Testing", hours filled in ----------------------
$ws.Range("A3").Value = "Testing"
$ws.Range("T3").Value = 4
$ws.Range("U3").Value = 8
$ws.Range("Y3").Value = 8
$ws.Range("Z3").Value = 5
$ws.Range("AA3").Value = 4
$ws.Range("AB3").Value = 4

# --- Row 4: hours filled in (task unchanged) --------------------------------
$ws.Range("Z4").Value = 3
$ws.Range("AA4").Value = 4
$ws.Range("AB4").Value = 4

# --- Column AF (last day of the month) gets the weekend/grey formatting ----
# used by the other non-entry columns, on every task row.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("AF2").PasteSpecial(-4122) | Out-Null

$ws.Range("B3").Copy() | Out-Null
$ws.Range("AF3").PasteSpecial(-4122) | Out-Null

$ws.Range("B4").Copy() | Out-Null
$ws.Range("AF4").PasteSpecial(-4122) | Out-Null

$ws.Range("B5").Copy() | Out-Null
$ws.Range("AF5").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Move the active selection on the visible sheet --------------------------
$ws.Range("T12").Select() | Out-Null
